# The title on slide 2 originally holds three runs ("Below", " ",
# "section-level") that together read "Below section-level". The target
# edit merges them into a single run with the same text. Because the
# COM text-setter preserves existing runs when the new text still shares
# a matching prefix/suffix with the old concatenation, just assigning the
# same final string back is a no-op for run layout. So we first stamp an
# unrelated placeholder (forcing the whole paragraph to collapse to one
# fresh run) and only then set the real text, which yields a single run.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange
$tr.Text = "x"
$tr.Text = "Below section-level"
